# chore: update Sheets via scheduled runner
#
# Refreshes the market-price / profit columns (H..N: currentAveragePrice,
# currentAveragePriceNQ, currentAveragePriceHQ, LevePriceNQ, LevePriceHQ,
# LeveProfitNQ, LeveProfitHQ) for the affected Leve rows across all eight
# job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly
# scraped market-board data. A few rows gain or lose an N (LeveProfitHQ)
# or M (LeveProfitNQ) cell entirely depending on whether an HQ/NQ price
# exists for that pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value2 = 428.8889  # ALC!H2: 458 -> 428.8889
$ws.Cells.Item(2, 9).Value2 = 431.33334  # ALC!I2: 458 -> 431.33334
$ws.Cells.Item(2, 10).Value2 = 424  # ALC!J2: 0 -> 424
$ws.Cells.Item(2, 11).Value2 = 431.33334  # ALC!K2: 458 -> 431.33334
$ws.Cells.Item(2, 12).Value2 = 424  # ALC!L2: 0 -> 424
$ws.Cells.Item(2, 13).Value2 = -318.33334  # ALC!M2: -345 -> -318.33334
$ws.Cells.Item(2, 14).Value2 = -650  # ALC!N2: <MISSING> -> -650
$ws.Cells.Item(9, 8).Value2 = 197.16667  # ALC!H9: 350 -> 197.16667
$ws.Cells.Item(9, 9).Value2 = 145.75  # ALC!I9: 200 -> 145.75
$ws.Cells.Item(9, 10).Value2 = 300  # ALC!J9: 500 -> 300
$ws.Cells.Item(9, 11).Value2 = 145.75  # ALC!K9: 200 -> 145.75
$ws.Cells.Item(9, 12).Value2 = 300  # ALC!L9: 500 -> 300
$ws.Cells.Item(9, 13).Value2 = 23.25  # ALC!M9: -31 -> 23.25
$ws.Cells.Item(9, 14).Value2 = -638  # ALC!N9: -838 -> -638
$ws.Cells.Item(28, 8).Value2 = 309.66666  # ALC!H28: 321.27274 -> 309.66666
$ws.Cells.Item(28, 9).Value2 = 342.125  # ALC!I28: 326 -> 342.125
$ws.Cells.Item(28, 10).Value2 = 244.75  # ALC!J28: 300 -> 244.75
$ws.Cells.Item(28, 11).Value2 = 342.125  # ALC!K28: 326 -> 342.125
$ws.Cells.Item(28, 12).Value2 = 244.75  # ALC!L28: 300 -> 244.75
$ws.Cells.Item(28, 13).Value2 = 142.875  # ALC!M28: 159 -> 142.875
$ws.Cells.Item(28, 14).Value2 = -1214.75  # ALC!N28: -1270 -> -1214.75
$ws.Cells.Item(38, 8).Value2 = 178.125  # ALC!H38: 301.5 -> 178.125
$ws.Cells.Item(38, 9).Value2 = 178.125  # ALC!I38: 301.5 -> 178.125
$ws.Cells.Item(38, 11).Value2 = 534.375  # ALC!K38: 904.5 -> 534.375
$ws.Cells.Item(38, 13).Value2 = -162.375  # ALC!M38: -532.5 -> -162.375
$ws.Cells.Item(40, 8).Value2 = 3250  # ALC!H40: 2624.875 -> 3250
$ws.Cells.Item(40, 9).Value2 = 2000  # ALC!I40: 1999.8 -> 2000
$ws.Cells.Item(40, 11).Value2 = 2000  # ALC!K40: 1999.8 -> 2000
$ws.Cells.Item(40, 13).Value2 = -1825  # ALC!M40: -1824.8 -> -1825
$ws.Cells.Item(58, 8).Value2 = 4382.2666  # ALC!H58: 4811.25 -> 4382.2666
$ws.Cells.Item(58, 9).Value2 = 2017.5  # ALC!I58: 35 -> 2017.5
$ws.Cells.Item(58, 10).Value2 = 4746.077  # ALC!J58: 5245.4546 -> 4746.077
$ws.Cells.Item(58, 11).Value2 = 6052.5  # ALC!K58: 105 -> 6052.5
$ws.Cells.Item(58, 12).Value2 = 14238.231  # ALC!L58: 15736.3638 -> 14238.231
$ws.Cells.Item(58, 13).Value2 = -5902.5  # ALC!M58: 45 -> -5902.5
$ws.Cells.Item(58, 14).Value2 = -14538.231  # ALC!N58: -16036.3638 -> -14538.231
$ws.Cells.Item(92, 8).Value2 = 959.6429000000001  # ALC!H92: 970.4286 -> 959.6429000000001
$ws.Cells.Item(92, 9).Value2 = 951.5  # ALC!I92: 1046.2222 -> 951.5
$ws.Cells.Item(92, 10).Value2 = 980  # ALC!J92: 834 -> 980
$ws.Cells.Item(92, 11).Value2 = 951.5  # ALC!K92: 1046.2222 -> 951.5
$ws.Cells.Item(92, 12).Value2 = 980  # ALC!L92: 834 -> 980
$ws.Cells.Item(92, 13).Value2 = 296.5  # ALC!M92: 201.7778000000001 -> 296.5
$ws.Cells.Item(92, 14).Value2 = -3476  # ALC!N92: -3330 -> -3476
$ws.Cells.Item(98, 8).Value2 = 1173.5454  # ALC!H98: 1290 -> 1173.5454
$ws.Cells.Item(98, 9).Value2 = 990.9  # ALC!I98: 1076.25 -> 990.9
$ws.Cells.Item(98, 11).Value2 = 990.9  # ALC!K98: 1076.25 -> 990.9
$ws.Cells.Item(98, 13).Value2 = 507.1  # ALC!M98: 421.75 -> 507.1
$ws.Cells.Item(122, 8).Value2 = 1173.5454  # ALC!H122: 1290 -> 1173.5454
$ws.Cells.Item(122, 9).Value2 = 990.9  # ALC!I122: 1076.25 -> 990.9
$ws.Cells.Item(122, 11).Value2 = 2972.7  # ALC!K122: 3228.75 -> 2972.7
$ws.Cells.Item(122, 13).Value2 = -522.6999999999998  # ALC!M122: -778.75 -> -522.6999999999998
$ws.Cells.Item(131, 8).Value2 = 4030.6667  # ALC!H131: 4090.8 -> 4030.6667
$ws.Cells.Item(131, 9).Value2 = 1079.4445  # ALC!I131: 1179.6666 -> 1079.4445
$ws.Cells.Item(131, 11).Value2 = 3238.3335  # ALC!K131: 3538.9998 -> 3238.3335
$ws.Cells.Item(131, 13).Value2 = 1801.6665  # ALC!M131: 1501.0002 -> 1801.6665
$ws.Cells.Item(137, 8).Value2 = 3012.5715  # ALC!H137: 3182.1667 -> 3012.5715
$ws.Cells.Item(137, 10).Value2 = 3012.5715  # ALC!J137: 3182.1667 -> 3012.5715
$ws.Cells.Item(137, 12).Value2 = 9037.7145  # ALC!L137: 9546.500100000001 -> 9037.7145
$ws.Cells.Item(137, 14).Value2 = -14137.7145  # ALC!N137: -14646.5001 -> -14137.7145
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 674.2  # ARM!H2: 691.3570999999999 -> 674.2
$ws.Cells.Item(2, 9).Value2 = 554.8461  # ARM!I2: 564.9167 -> 554.8461
$ws.Cells.Item(2, 11).Value2 = 554.8461  # ARM!K2: 564.9167 -> 554.8461
$ws.Cells.Item(2, 13).Value2 = -441.8461  # ARM!M2: -451.9167 -> -441.8461
$ws.Cells.Item(45, 8).Value2 = 7635.6665  # ARM!H45: 4767.3335 -> 7635.6665
$ws.Cells.Item(45, 9).Value2 = 2204.5  # ARM!I45: 2021.2 -> 2204.5
$ws.Cells.Item(45, 11).Value2 = 2204.5  # ARM!K45: 2021.2 -> 2204.5
$ws.Cells.Item(45, 13).Value2 = -1827.5  # ARM!M45: -1644.2 -> -1827.5
$ws.Cells.Item(116, 8).Value2 = 674.2  # ARM!H116: 691.3570999999999 -> 674.2
$ws.Cells.Item(116, 9).Value2 = 554.8461  # ARM!I116: 564.9167 -> 554.8461
$ws.Cells.Item(116, 11).Value2 = 554.8461  # ARM!K116: 564.9167 -> 554.8461
$ws.Cells.Item(116, 13).Value2 = 1739.1539  # ARM!M116: 1729.0833 -> 1739.1539
$ws.Cells.Item(132, 8).Value2 = 1307.7441  # ARM!H132: 1254.4565 -> 1307.7441
$ws.Cells.Item(132, 9).Value2 = 1315.0952  # ARM!I132: 1295.8372 -> 1315.0952
$ws.Cells.Item(132, 10).Value2 = 999  # ARM!J132: 661.3333 -> 999
$ws.Cells.Item(132, 11).Value2 = 3945.2856  # ARM!K132: 3887.5116 -> 3945.2856
$ws.Cells.Item(132, 12).Value2 = 2997  # ARM!L132: 1983.9999 -> 2997
$ws.Cells.Item(132, 13).Value2 = -1415.2856  # ARM!M132: -1357.5116 -> -1415.2856
$ws.Cells.Item(132, 14).Value2 = -8057  # ARM!N132: -7043.9999 -> -8057
$ws.Cells.Item(135, 8).Value2 = 0  # ARM!H135: 47500 -> 0
$ws.Cells.Item(135, 10).Value2 = 0  # ARM!J135: 47500 -> 0
$ws.Cells.Item(135, 12).Value2 = 0  # ARM!L135: 47500 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # ARM!N135: -57640 -> (removed)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 674.2  # BSM!H3: 691.3570999999999 -> 674.2
$ws.Cells.Item(3, 9).Value2 = 554.8461  # BSM!I3: 564.9167 -> 554.8461
$ws.Cells.Item(3, 11).Value2 = 554.8461  # BSM!K3: 564.9167 -> 554.8461
$ws.Cells.Item(3, 13).Value2 = -440.8461  # BSM!M3: -450.9167 -> -440.8461
$ws.Cells.Item(5, 8).Value2 = 1749.75  # BSM!H5: 1579.8 -> 1749.75
$ws.Cells.Item(5, 9).Value2 = 1999.6666  # BSM!I5: 1724.75 -> 1999.6666
$ws.Cells.Item(5, 11).Value2 = 1999.6666  # BSM!K5: 1724.75 -> 1999.6666
$ws.Cells.Item(5, 13).Value2 = -1886.6666  # BSM!M5: -1611.75 -> -1886.6666
$ws.Cells.Item(94, 8).Value2 = 1197.6  # BSM!H94: 1199 -> 1197.6
$ws.Cells.Item(94, 9).Value2 = 999.3333  # BSM!I94: 998.75 -> 999.3333
$ws.Cells.Item(94, 10).Value2 = 1495  # BSM!J94: 2000 -> 1495
$ws.Cells.Item(94, 11).Value2 = 999.3333  # BSM!K94: 998.75 -> 999.3333
$ws.Cells.Item(94, 12).Value2 = 1495  # BSM!L94: 2000 -> 1495
$ws.Cells.Item(94, 13).Value2 = -548.3333  # BSM!M94: -547.75 -> -548.3333
$ws.Cells.Item(94, 14).Value2 = -2397  # BSM!N94: -2902 -> -2397
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 4272.391  # CRP!H31: 4397.591 -> 4272.391
$ws.Cells.Item(31, 9).Value2 = 2780.7778  # CRP!I31: 2834.3333 -> 2780.7778
$ws.Cells.Item(31, 10).Value2 = 5231.2856  # CRP!J31: 5479.846 -> 5231.2856
$ws.Cells.Item(31, 11).Value2 = 2780.7778  # CRP!K31: 2834.3333 -> 2780.7778
$ws.Cells.Item(31, 12).Value2 = 5231.2856  # CRP!L31: 5479.846 -> 5231.2856
$ws.Cells.Item(31, 13).Value2 = -2485.7778  # CRP!M31: -2539.3333 -> -2485.7778
$ws.Cells.Item(31, 14).Value2 = -5821.2856  # CRP!N31: -6069.846 -> -5821.2856
$ws.Cells.Item(34, 8).Value2 = 4272.391  # CRP!H34: 4397.591 -> 4272.391
$ws.Cells.Item(34, 9).Value2 = 2780.7778  # CRP!I34: 2834.3333 -> 2780.7778
$ws.Cells.Item(34, 10).Value2 = 5231.2856  # CRP!J34: 5479.846 -> 5231.2856
$ws.Cells.Item(34, 11).Value2 = 2780.7778  # CRP!K34: 2834.3333 -> 2780.7778
$ws.Cells.Item(34, 12).Value2 = 5231.2856  # CRP!L34: 5479.846 -> 5231.2856
$ws.Cells.Item(34, 13).Value2 = -2578.7778  # CRP!M34: -2632.3333 -> -2578.7778
$ws.Cells.Item(34, 14).Value2 = -5635.2856  # CRP!N34: -5883.846 -> -5635.2856
$ws.Cells.Item(44, 8).Value2 = 16000  # CRP!H44: 11021.667 -> 16000
$ws.Cells.Item(44, 9).Value2 = 15000  # CRP!I44: 65 -> 15000
$ws.Cells.Item(44, 11).Value2 = 15000  # CRP!K44: 65 -> 15000
$ws.Cells.Item(44, 13).Value2 = -14558  # CRP!M44: 377 -> -14558
$ws.Cells.Item(68, 8).Value2 = 28374.75  # CRP!H68: 31999.666 -> 28374.75
$ws.Cells.Item(68, 9).Value2 = 23500  # CRP!I68: 0 -> 23500
$ws.Cells.Item(68, 10).Value2 = 29999.666  # CRP!J68: 31999.666 -> 29999.666
$ws.Cells.Item(68, 11).Value2 = 23500  # CRP!K68: 0 -> 23500
$ws.Cells.Item(68, 12).Value2 = 29999.666  # CRP!L68: 31999.666 -> 29999.666
$ws.Cells.Item(68, 13).Value2 = -22751  # CRP!M68: <MISSING> -> -22751
$ws.Cells.Item(68, 14).Value2 = -31497.666  # CRP!N68: -33497.666 -> -31497.666
$ws.Cells.Item(71, 8).Value2 = 28374.75  # CRP!H71: 31999.666 -> 28374.75
$ws.Cells.Item(71, 9).Value2 = 23500  # CRP!I71: 0 -> 23500
$ws.Cells.Item(71, 10).Value2 = 29999.666  # CRP!J71: 31999.666 -> 29999.666
$ws.Cells.Item(71, 11).Value2 = 70500  # CRP!K71: 0 -> 70500
$ws.Cells.Item(71, 12).Value2 = 89998.99800000001  # CRP!L71: 95998.99800000001 -> 89998.99800000001
$ws.Cells.Item(71, 13).Value2 = -66756  # CRP!M71: <MISSING> -> -66756
$ws.Cells.Item(71, 14).Value2 = -97486.99800000001  # CRP!N71: -103486.998 -> -97486.99800000001
$ws.Cells.Item(134, 8).Value2 = 4296.5  # CRP!H134: 4849.5 -> 4296.5
$ws.Cells.Item(134, 9).Value2 = 2882.8  # CRP!I134: 3415 -> 2882.8
$ws.Cells.Item(134, 11).Value2 = 8648.400000000001  # CRP!K134: 10245 -> 8648.400000000001
$ws.Cells.Item(134, 13).Value2 = -6113.400000000001  # CRP!M134: -7710 -> -6113.400000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value2 = 6666738.5  # CUL!H7: 7692390 -> 6666738.5
$ws.Cells.Item(7, 9).Value2 = 9090962  # CUL!I7: 10000057 -> 9090962
$ws.Cells.Item(7, 10).Value2 = 125  # CUL!J7: 166.33333 -> 125
$ws.Cells.Item(7, 11).Value2 = 27272886  # CUL!K7: 30000171 -> 27272886
$ws.Cells.Item(7, 12).Value2 = 375  # CUL!L7: 498.99999 -> 375
$ws.Cells.Item(7, 13).Value2 = -27272774  # CUL!M7: -30000059 -> -27272774
$ws.Cells.Item(7, 14).Value2 = -599  # CUL!N7: -722.99999 -> -599
$ws.Cells.Item(75, 8).Value2 = 350.625  # CUL!H75: 343.57144 -> 350.625
$ws.Cells.Item(75, 9).Value2 = 392.5  # CUL!I75: 391 -> 392.5
$ws.Cells.Item(75, 11).Value2 = 1177.5  # CUL!K75: 1173 -> 1177.5
$ws.Cells.Item(75, 13).Value2 = -179.5  # CUL!M75: -175 -> -179.5
$ws.Cells.Item(78, 8).Value2 = 350.625  # CUL!H78: 343.57144 -> 350.625
$ws.Cells.Item(78, 9).Value2 = 392.5  # CUL!I78: 391 -> 392.5
$ws.Cells.Item(78, 11).Value2 = 3532.5  # CUL!K78: 3519 -> 3532.5
$ws.Cells.Item(78, 13).Value2 = 1459.5  # CUL!M78: 1473 -> 1459.5
$ws.Cells.Item(107, 8).Value2 = 1262.1  # CUL!H107: 1136.75 -> 1262.1
$ws.Cells.Item(107, 9).Value2 = 669  # CUL!I107: 589.5 -> 669
$ws.Cells.Item(107, 11).Value2 = 2007  # CUL!K107: 1768.5 -> 2007
$ws.Cells.Item(107, 13).Value2 = -87  # CUL!M107: 151.5 -> -87
$ws.Cells.Item(129, 8).Value2 = 1949.2727  # CUL!H129: 2629.7856 -> 1949.2727
$ws.Cells.Item(129, 10).Value2 = 1879  # CUL!J129: 3096.25 -> 1879
$ws.Cells.Item(129, 12).Value2 = 5637  # CUL!L129: 9288.75 -> 5637
$ws.Cells.Item(129, 14).Value2 = -15637  # CUL!N129: -19288.75 -> -15637
$ws.Cells.Item(131, 8).Value2 = 3701.2666  # CUL!H131: 3554 -> 3701.2666
$ws.Cells.Item(131, 9).Value2 = 3231.75  # CUL!I131: 3022.111 -> 3231.75
$ws.Cells.Item(131, 11).Value2 = 9695.25  # CUL!K131: 9066.332999999999 -> 9695.25
$ws.Cells.Item(131, 13).Value2 = -4655.25  # CUL!M131: -4026.332999999999 -> -4655.25
$ws.Cells.Item(132, 8).Value2 = 8018.2  # CUL!H132: 4322 -> 8018.2
$ws.Cells.Item(132, 9).Value2 = 0  # CUL!I132: 634 -> 0
$ws.Cells.Item(132, 10).Value2 = 8018.2  # CUL!J132: 6780.6665 -> 8018.2
$ws.Cells.Item(132, 11).Value2 = 0  # CUL!K132: 5706 -> 0
$ws.Cells.Item(132, 12).Value2 = 72163.8  # CUL!L132: 61025.9985 -> 72163.8
$ws.Cells.Item(132, 13).ClearContents()  # CUL!M132: -3176 -> (removed)
$ws.Cells.Item(132, 14).Value2 = -77223.8  # CUL!N132: -66085.9985 -> -77223.8
$ws.Cells.Item(139, 8).Value2 = 3052.1333  # CUL!H139: 2891.6428 -> 3052.1333
$ws.Cells.Item(139, 9).Value2 = 3052.1333  # CUL!I139: 2891.6428 -> 3052.1333
$ws.Cells.Item(139, 11).Value2 = 9156.3999  # CUL!K139: 8674.928400000001 -> 9156.3999
$ws.Cells.Item(139, 13).Value2 = -4016.3999  # CUL!M139: -3534.928400000001 -> -4016.3999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value2 = 2626.0625  # GSM!H113: 2457.111 -> 2626.0625
$ws.Cells.Item(113, 9).Value2 = 1201.7  # GSM!I113: 1185.6666 -> 1201.7
$ws.Cells.Item(113, 11).Value2 = 1201.7  # GSM!K113: 1185.6666 -> 1201.7
$ws.Cells.Item(113, 13).Value2 = 968.3  # GSM!M113: 984.3334 -> 968.3
$ws.Cells.Item(122, 8).Value2 = 690272  # GSM!H122: 614013.9399999999 -> 690272
$ws.Cells.Item(122, 9).Value2 = 127312.25  # GSM!I122: 102639.7 -> 127312.25
$ws.Cells.Item(122, 11).Value2 = 381936.75  # GSM!K122: 307919.1 -> 381936.75
$ws.Cells.Item(122, 13).Value2 = -379486.75  # GSM!M122: -305469.1 -> -379486.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 967.1667  # LTW!H22: 1501 -> 967.1667
$ws.Cells.Item(22, 9).Value2 = 1050.25  # LTW!I22: 2000 -> 1050.25
$ws.Cells.Item(22, 10).Value2 = 801  # LTW!J22: 1002 -> 801
$ws.Cells.Item(22, 11).Value2 = 1050.25  # LTW!K22: 2000 -> 1050.25
$ws.Cells.Item(22, 12).Value2 = 801  # LTW!L22: 1002 -> 801
$ws.Cells.Item(22, 13).Value2 = -755.25  # LTW!M22: -1705 -> -755.25
$ws.Cells.Item(22, 14).Value2 = -1391  # LTW!N22: -1592 -> -1391
$ws.Cells.Item(27, 8).Value2 = 967.1667  # LTW!H27: 1501 -> 967.1667
$ws.Cells.Item(27, 9).Value2 = 1050.25  # LTW!I27: 2000 -> 1050.25
$ws.Cells.Item(27, 10).Value2 = 801  # LTW!J27: 1002 -> 801
$ws.Cells.Item(27, 11).Value2 = 1050.25  # LTW!K27: 2000 -> 1050.25
$ws.Cells.Item(27, 12).Value2 = 801  # LTW!L27: 1002 -> 801
$ws.Cells.Item(27, 13).Value2 = -943.25  # LTW!M27: -1893 -> -943.25
$ws.Cells.Item(27, 14).Value2 = -1015  # LTW!N27: -1216 -> -1015
$ws.Cells.Item(61, 8).Value2 = 2227.5  # LTW!H61: 2296.0667 -> 2227.5
$ws.Cells.Item(61, 9).Value2 = 1553.1666  # LTW!I61: 1585.3636 -> 1553.1666
$ws.Cells.Item(61, 11).Value2 = 1553.1666  # LTW!K61: 1585.3636 -> 1553.1666
$ws.Cells.Item(61, 13).Value2 = -1351.1666  # LTW!M61: -1383.3636 -> -1351.1666
$ws.Cells.Item(82, 8).Value2 = 2185.2666  # LTW!H82: 2277.4285 -> 2185.2666
$ws.Cells.Item(82, 10).Value2 = 2199  # LTW!J82: 2525 -> 2199
$ws.Cells.Item(82, 12).Value2 = 2199  # LTW!L82: 2525 -> 2199
$ws.Cells.Item(82, 14).Value2 = -2921  # LTW!N82: -3247 -> -2921
$ws.Cells.Item(85, 8).Value2 = 2185.2666  # LTW!H85: 2277.4285 -> 2185.2666
$ws.Cells.Item(85, 10).Value2 = 2199  # LTW!J85: 2525 -> 2199
$ws.Cells.Item(85, 12).Value2 = 2199  # LTW!L85: 2525 -> 2199
$ws.Cells.Item(85, 14).Value2 = -4695  # LTW!N85: -5021 -> -4695
$ws.Cells.Item(113, 8).Value2 = 2227.5  # LTW!H113: 2296.0667 -> 2227.5
$ws.Cells.Item(113, 9).Value2 = 1553.1666  # LTW!I113: 1585.3636 -> 1553.1666
$ws.Cells.Item(113, 11).Value2 = 1553.1666  # LTW!K113: 1585.3636 -> 1553.1666
$ws.Cells.Item(113, 13).Value2 = 616.8334  # LTW!M113: 584.6364000000001 -> 616.8334
$ws.Cells.Item(132, 8).Value2 = 4706.5117  # LTW!H132: 4787.884 -> 4706.5117
$ws.Cells.Item(132, 9).Value2 = 4349.8965  # LTW!I132: 4470.552 -> 4349.8965
$ws.Cells.Item(132, 11).Value2 = 13049.6895  # LTW!K132: 13411.656 -> 13049.6895
$ws.Cells.Item(132, 13).Value2 = -10519.6895  # LTW!M132: -10881.656 -> -10519.6895
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(60, 8).Value2 = 99997  # WVR!H60: 99997.664 -> 99997
$ws.Cells.Item(60, 10).Value2 = 100000  # WVR!J60: 99999.5 -> 100000
$ws.Cells.Item(60, 12).Value2 = 100000  # WVR!L60: 99999.5 -> 100000
$ws.Cells.Item(60, 14).Value2 = -101644  # WVR!N60: -101643.5 -> -101644
$ws.Cells.Item(62, 8).Value2 = 8155.6875  # WVR!H62: 7624.375 -> 8155.6875
$ws.Cells.Item(62, 10).Value2 = 8155.6875  # WVR!J62: 7624.375 -> 8155.6875
$ws.Cells.Item(62, 12).Value2 = 8155.6875  # WVR!L62: 7624.375 -> 8155.6875
$ws.Cells.Item(62, 14).Value2 = -9403.6875  # WVR!N62: -8872.375 -> -9403.6875
$ws.Cells.Item(65, 8).Value2 = 8155.6875  # WVR!H65: 7624.375 -> 8155.6875
$ws.Cells.Item(65, 10).Value2 = 8155.6875  # WVR!J65: 7624.375 -> 8155.6875
$ws.Cells.Item(65, 12).Value2 = 40778.4375  # WVR!L65: 38121.875 -> 40778.4375
$ws.Cells.Item(65, 14).Value2 = -47018.4375  # WVR!N65: -44361.875 -> -47018.4375
$ws.Cells.Item(96, 8).Value2 = 1639.5  # WVR!H96: 1636.5555 -> 1639.5
$ws.Cells.Item(96, 9).Value2 = 1619.1666  # WVR!I96: 1609.8 -> 1619.1666
$ws.Cells.Item(96, 11).Value2 = 1619.1666  # WVR!K96: 1609.8 -> 1619.1666
$ws.Cells.Item(96, 13).Value2 = -246.1666  # WVR!M96: -236.8 -> -246.1666
